$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row for "디비금융스팩12호" (DB Financial SPAC 12) above the
# existing row 3 ("그리드위즈"), shifting all subsequent IPO rows down by one.
$ws.Rows.Item(3).Insert()

# Date-like text must stay plain text (shared string), not be auto-converted
# to a date serial number by the smart-parsing Value setter.
$ws.Range("A3:E3").NumberFormat = "@"

$ws.Range("A3").Value = "2024-06-05"
$ws.Range("B3").Value = "디비금융스팩12호"
$ws.Range("C3").Value = "DB"
$ws.Range("D3").Value = "2024-06-11"
$ws.Range("E3").Value = "2024-06-18"
$ws.Range("F3").Value = 10000000
$ws.Range("G3").Value = 5000000
$ws.Range("H3").Value = "-"
$ws.Range("I3").Value = 2000
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = "-"
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = "-"
$ws.Range("N3").Value = "-"
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = "-"
$ws.Range("Q3").Value = "-"
$ws.Range("R3").Value = "1141.4 : 1"
$ws.Range("S3").Value = "-"
$ws.Range("T3").Value = "-"

# Drop the temporary text-format override so the new row matches the plain
# (unstyled) look of every other data row.
$ws.Range("A3:T3").Style = "Normal"
